$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text fixes (affects every cell sharing the old string) ---

# "<1-01-31" -> "<001-01-31"  (F2, F3)
$ws.Range("F2").Value = "<001-01-31"
$ws.Range("F3").Value = "<001-01-31"

# "z0bug.jou_inv" -> "external.FAT|FATT|INV"  (K2:K9, K11)
$ws.Range("K2").Value = "external.FAT|FATT|INV"
$ws.Range("K3").Value = "external.FAT|FATT|INV"
$ws.Range("K4").Value = "external.FAT|FATT|INV"
$ws.Range("K5").Value = "external.FAT|FATT|INV"
$ws.Range("K6").Value = "external.FAT|FATT|INV"
$ws.Range("K7").Value = "external.FAT|FATT|INV"
$ws.Range("K8").Value = "external.FAT|FATT|INV"
$ws.Range("K9").Value = "external.FAT|FATT|INV"
$ws.Range("K11").Value = "external.FAT|FATT|INV"

# "<2-12-20" -> "<002-12-20"  (F12)
$ws.Range("F12").Value = "<002-12-20"

# "<2-12-99" -> "<002-12-99"  (G12, I12, G13)
$ws.Range("G12").Value = "<002-12-99"
$ws.Range("I12").Value = "<002-12-99"
$ws.Range("G13").Value = "<002-12-99"

# "<1-01-01" -> "<001-01-01"  (H12, H13)
$ws.Range("H12").Value = "<001-01-01"
$ws.Range("H13").Value = "<001-01-01"

# "z0bug.jou_bill" -> "external.ACQ|FATTU|BILL"  (K12:K18, K20)
$ws.Range("K12").Value = "external.ACQ|FATTU|BILL"
$ws.Range("K13").Value = "external.ACQ|FATTU|BILL"
$ws.Range("K14").Value = "external.ACQ|FATTU|BILL"
$ws.Range("K15").Value = "external.ACQ|FATTU|BILL"
$ws.Range("K16").Value = "external.ACQ|FATTU|BILL"
$ws.Range("K17").Value = "external.ACQ|FATTU|BILL"
$ws.Range("K18").Value = "external.ACQ|FATTU|BILL"
$ws.Range("K20").Value = "external.ACQ|FATTU|BILL"

# "<2-12-22" -> "<002-12-22"  (F13)
$ws.Range("F13").Value = "<002-12-22"

# --- Style changes ---

# journal_id column (K) cells that now hold the new external codes pick up
# the "colored arial 9" style (s=2 -> s=5, i.e. font color turns explicit black)
$ws.Range("K2:K9").Font.Color = 0
$ws.Range("K11:K18").Font.Color = 0
$ws.Range("K20").Font.Color = 0

# payment_term_id column (O) cells move from the arial-9 style to the new
# Calibri-9 style (s=5 -> s=6, i.e. a new font is introduced)
$ws.Range("O3").Font.Name = "Calibri"
$ws.Range("O5").Font.Name = "Calibri"
$ws.Range("O7").Font.Name = "Calibri"
$ws.Range("O9").Font.Name = "Calibri"
$ws.Range("O12:O20").Font.Name = "Calibri"

# --- Column width (K / column 11) ---
$ws.Columns.Item(11).ColumnWidth = 16.8

# --- Selection / frozen-pane bottom-right anchor ---
$ws.Range("K20").Select()
